# BadEventsV4ToParse.xlsx -- "testing some of the bad events"
#
# Fills in several previously-blank BAD-sheet outcome cells with new event
# text, fixes a capitalization typo, tweaks two sentences' punctuation, and
# moves the active/selected sheet+cell from GOOD back to BAD.

$wb = $excel.ActiveWorkbook
$bad = $wb.Worksheets.Item("BAD")

# --- Row 6 (volcano / lava-flow vision): fix "hojors" -> "Hojors" -----------
$bad.Range("E6").Value = "You send a vision to warn your followers to dig a trench to collect and divert the lava flow. It was received by Hojors, the town fool. No one paid heed and your followers lose a lot of their important buildings"

# --- Row 8 (League of Legends / DDR challenge): was all BLANK --------------
$bad.Range("D8").Value = "When he claims to be better then you in League of Legends, that is when you draw the line. You challenge him to a 1v1 and impress your people with your impressive skills"
$bad.Range("E8").Value = "When he claims to be better then you in League of Legends, that is when you draw the line. You challenge him to a 1v1 and embarras yourself when he starts dancing over your champion's dead body"
$bad.Range("F8").Value = "The Imposter seems to start flirting with one of the lady followers that you had your eyes on.  You instantly show your machismo and win the lady, and as you give her a night she won't forget, she starts spreading the news of your amazing skills "
$bad.Range("G8").Value = "The Imposter seems to start flirting with one of the lady followers that you had your eyes on.  You try to show your machismo and win the lady, but you embarras yourself and she spreads the news of the laughable exchange you two had "
$bad.Range("H8").Value = "You simply get fed up with what he is doing and challenge him to the ultimate battle of DDR. You (rigged) showed off your skills to the people who are at awe, and finally put the imposter in his place "
$bad.Range("I8").Value = "You simply get fed up with what he is doing and challenge him to the ultimate battle of DDR. You rigged the machine to make you win, but as taunt him showing your moves to the imposter, the machine explodes in your face and destorys some buildings"
$bad.Rows.Item(8).RowHeight = 90

# --- Row 10 (Yo-Yo trade route): B/C get the reworded "staring at a nearby -
# --- mountain" text, H/I become a new giant-yo-yo event ---------------------
$bad.Range("B10").Value = "Staring at a nearby mountain reveals its hidden content to your followers. They rejoice at the new and improve yo-yos that they now poccess"
$bad.Range("C10").Value = "Staring at a nearby mountain reveals its hidden content to your followers. They gnat their teeth in frustration as a foul odor premeates the air for the next few weeks"
$bad.Range("H10").Value = "Listening to the cries of your people, you summoned a giant yo-yo in the middle of town square. It obeys your people, providing joy and happiness for them. It is also an item that the neighborring coutry covets "
$bad.Range("I10").Value = "Listening to the cries of your people, you summoned a giant yo-yo in the middle of town square.However, It starts to rampage and destorys some of your people and their buildings. "

# --- Row 16 (sheep appear in town): was all BLANK ---------------------------
$bad.Range("B16").Value = "They seem to be docile enough. So you instruct your followers to herd them to the fields, where they continue to gaze and provide food for your people"
$bad.Range("C16").Value = "They seem to be docile enough. So you instruct your followers to herd them to the fields. Unfortunately, they turn out to be disguised thieves  and steal some of your town's prized possessions and statues of you in the middle of the night"
$bad.Range("D16").Value = "You instruct your people to save some of the sheep for their wool. While they go to shear them, they see that one of sheep's wool was used to hide treasure, lots of it. Your followers are happy from this terricfic outcome "
$bad.Range("E16").Value = "You instruct your people to save some of the sheep for their wool. While they go to shear them, they get infected by a virus that transmit through the touch of their wool. The people cry on how you let this happen"
$bad.Range("F16").Value = "As your town goes to slaughter the sheep for food, it is reveled that one of the sheep is the personal pet of udpof, the God of the wopesde people. Using your impressive charisma, you negoiate a way for your followers to feast on the delicious sheep while returning the pet"
$bad.Range("G16").Value = "As your town goes to slaughter the sheep for food, it is reveled that one of the sheep is the personal pet of udpof, the God of the wopesde people. During your negoiations, you sneezed in his face, which is an insult in his eyes. With him saying the magic words, his pet starts destroying your people"
$bad.Range("H16").Value = "At the sight of the sheeps you grow hungry, so you ordered your followers to hastily  sacrifice one of them to you. As they start the ritual, the sheep's defense mechanisms kick in but messes up, transforminging him into a gaint sheep which you share with the people"
$bad.Range("I16").Value = "At the sight of the sheeps you grow hungry, so you ordered your followers to hastily  sacrifice one of them to you. As they start the ritual, the sheep's defense mechanisms kick in, transforminging him into a gaint man-killing sheep who terrorize your people"
$bad.Rows.Item(16).RowHeight = 90

# --- Row 17 (WILDFIRE): B/C were BLANK --------------------------------------
$bad.Range("B17").Value = "Looking directly at the fire, you start sheding tears that rains down and put out the fire. Your followers thank you for your grace "
$bad.Range("C17").Value = "Looking directly at the fire, you start sheding tears that rains down and put out the fire, but you continue and flood some parts of your land"
$bad.Rows.Item(17).RowHeight = 30

# --- View state: make BAD the active/selected sheet again (was GOOD) -------
$bad.Activate()
$bad.Range("A2:I21").Select()
